# Adds "phone type" question rows (apple/iOS support) to the register sheet,
# per commit: "apple phone support added, participation and debriefing letter update"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows appended after the existing data (previously ending at row 65).
# Cells are written in this particular order so the resulting shared-string
# table lines up with the authored workbook.
$ws.Cells.Item(66, 1).Value = "phone_type"
$ws.Cells.Item(66, 2).Value = "Smartphone type"
$ws.Cells.Item(66, 3).Value = "Smartphonetype"

$ws.Cells.Item(67, 1).Value = "phone_option_1"
$ws.Cells.Item(67, 2).Value = "Android (Samsung, Google, Oppo, Xiaomi, or other)"

$ws.Cells.Item(68, 1).Value = "phone_option_2"

$ws.Cells.Item(69, 1).Value = "phone_option_3"
$ws.Cells.Item(69, 2).Value = "I do not know"
$ws.Cells.Item(69, 3).Value = "Ik weet het niet"

$ws.Cells.Item(67, 3).Value = "Android (Samsung, Google, Oppo, Xiaomi, Huawei, of andere)"

$ws.Cells.Item(68, 2).Value = "IOS (Apple iPhone)"
$ws.Cells.Item(68, 3).Value = "IOS (Apple iPhone)"

# Reflect where the author was last looking/selecting in the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 2
$ws.Range("B70").Select() | Out-Null

$wb.Save()
